$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = -22.06000000000001
$ws.Range("B4").Value = 5.133100000000001
$ws.Range("A7").Value = -20.11619999999998
$ws.Range("B12").Value = 4.8659
$ws.Range("A16").Value = -22.01100000000002
$ws.Range("B18").Value = 6.082999999999995
$ws.Range("B19").Value = 8.563900000000006
$ws.Range("B20").Value = 9.535099999999989
$ws.Range("A28").Value = -21.8009
$ws.Range("A29").Value = -21.09449999999996
$ws.Range("B31").Value = 5.514799999999999
$ws.Range("A32").Value = -21.20500000000001
$ws.Range("A40").Value = -20.69959999999997
$ws.Range("B40").Value = 7.915899999999998
$ws.Range("B42").Value = 8.638100000000003
$ws.Range("B47").Value = 5.929600000000004
$ws.Range("B48").Value = 5.844600000000003
$ws.Range("A52").Value = -22.2722
$ws.Range("A57").Value = -22.16370000000002
$ws.Range("B63").Value = 4.852799999999998
$ws.Range("B64").Value = 5.324000000000001
$ws.Range("A66").Value = -21.45310000000001
$ws.Range("B76").Value = 5.567999999999997
$ws.Range("B81").Value = 5.370800000000006
$ws.Range("B89").Value = 4.554699999999994
$ws.Range("B94").Value = 4.748699999999995
$ws.Range("A100").Value = -22.03460000000002
